# Update the last row of the tracking table ("Generar Informe Gerencial")
# to reflect the results of the initial test run, and mark every cell in
# that row with the green "done" shading (RGB 00B050).

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)
$row = $table.Rows.Item($table.Rows.Count)

# New cell texts, in column order (1-based). Column 1 (the CU name) keeps
# its existing text and only gets the new shading.
$newTexts = @{
    2  = "Realizada"
    3  = "N/A"
    4  = "15/10/2020"
    5  = "N/A"
    6  = "N/A"
    7  = "N/A"
    8  = "N/A"
    9  = "N/A"
    10 = "Aprobado (CU cerrado - No se debe modificar)"
}

# BGR-encoded 00B050 (Word COM colors are 0xBBGGRR) -> 0x50B000
$greenFill = 0x50B000

for ($i = 1; $i -le $row.Cells.Count; $i++) {
    $cell = $row.Cells.Item($i)
    $cell.Shading.BackgroundPatternColor = $greenFill
    if ($newTexts.ContainsKey($i)) {
        $cell.Range.Text = $newTexts[$i]
    }
}
